$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "Group 2 Diagnostics" row (smdi_rf group description):
# "Assessing the ability to predic missingness based on observed covariates"
#                                  ^ missing "t"          ^ narrow no-break spaces
# becomes
# "Assessing the ability to predict missingness based on observed covariates"
$ws.Range("C5").Value = "Assessing the ability to predict missingness based on observed covariates"

# Restore the view: scroll so column A is visible again and select C5
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 2
